# Applies the "Anonimyzed fedcore" update:
#  - renames the "fedcore" column header to "approach" on both sheets
#  - adds header/divider borders (top+bottom, and top+bottom+right) to the
#    small numeric cells that sit above the metric table header row
#  - clears the stray empty inline-string cell G5 on the computational sheet
#  - (the "-0" -> "0" numeric normalization happens automatically when the
#    workbook is re-serialized, no explicit action required)

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # quality_comparison
$ws2 = $wb.Worksheets.Item(2)   # computational_comparison

# ------------------------------------------------------------------
# 1. Build the two new border styles once (on sheet 1) and reuse them
#    everywhere else via copy/paste-special of formats only, so the
#    workbook ends up with exactly two new cell styles:
#      - top+bottom border                (used by C1 / F1)
#      - top+bottom+right border          (used by D1 / G1)
# ------------------------------------------------------------------

$topBottom = $ws1.Range("C1")
$topBottom.Style = "Normal"
$topBottom.Borders.Item(8).LineStyle = 1    # xlEdgeTop
$topBottom.Borders.Item(9).LineStyle = 1    # xlEdgeBottom

$topBottomRight = $ws1.Range("D1")
$topBottomRight.Style = "Normal"
$topBottomRight.Borders.Item(8).LineStyle = 1   # xlEdgeTop
$topBottomRight.Borders.Item(10).LineStyle = 1  # xlEdgeRight
$topBottomRight.Borders.Item(9).LineStyle = 1   # xlEdgeBottom

# Propagate the same two styles to the matching cells on sheet 2
$topBottom.Copy()
$ws2.Range("C1").PasteSpecial(-4122)   # xlPasteFormats
$ws2.Range("F1").PasteSpecial(-4122)

$topBottomRight.Copy()
$ws2.Range("D1").PasteSpecial(-4122)
$ws2.Range("G1").PasteSpecial(-4122)

# ------------------------------------------------------------------
# 2. Anonymize: rename "fedcore" -> "approach"
# ------------------------------------------------------------------
$ws1.Range("C2").Value2 = "approach"

$ws2.Range("C2").Value2 = "approach"
$ws2.Range("F2").Value2 = "approach"

# ------------------------------------------------------------------
# 3. Drop the stray empty cell G5 on the computational sheet
# ------------------------------------------------------------------
$ws2.Range("G5").ClearContents()
